$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row for id=5 (차두리): defense (F5) and stamina (G5) drop from 80 to 50
$ws.Range("F5").Value = 50.0
$ws.Range("G5").Value = 50.0

# Add a new player row: id=7, 박지성, with stats all 80.
# Copy the last data row down first so the new row inherits the same
# cell formatting/style as the rest of the table, then overwrite the values.
$ws.Rows("6:6").Copy()
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = 7.0
$ws.Range("B7").Value = "박지성"
$ws.Range("C7").Value = 80.0
$ws.Range("D7").Value = 80.0
$ws.Range("E7").Value = 80.0
$ws.Range("F7").Value = 80.0
$ws.Range("G7").Value = 80.0
